# Update PowerShell scripts used to create teams
# Rename the team entries in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Team 51"
$ws.Range("A3").Value = "Team 52"
$ws.Range("A4").Value = "Team 53"

$ws.Range("C6").Select()
